$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 57536.5
$ws.Range("J3").Value = 57536.5
$ws.Range("L3").Value = 57536.5
$ws.Range("N3").Value = -57764.5

$ws.Range("H17").Value = 445.25
$ws.Range("J17").Value = 449.7857
$ws.Range("L17").Value = 1349.3571
$ws.Range("N17").Value = -1685.3571

$ws.Range("H80").Value = 5604080
$ws.Range("J80").Value = 5954389.5
$ws.Range("L80").Value = 17863168.5
$ws.Range("N80").Value = -17865164.5

$ws.Range("H83").Value = 5604080
$ws.Range("J83").Value = 5954389.5
$ws.Range("L83").Value = 53589505.5
$ws.Range("N83").Value = -53599489.5

$ws.Range("H86").Value = 4344.4287
$ws.Range("I86").Value = 2475
$ws.Range("K86").Value = 2475
$ws.Range("M86").Value = -1352

$ws.Range("H88").Value = 2979.25
$ws.Range("J88").Value = 2945.875
$ws.Range("L88").Value = 2945.875
$ws.Range("N88").Value = -3757.875

$ws.Range("H89").Value = 4344.4287
$ws.Range("I89").Value = 2475
$ws.Range("K89").Value = 12375
$ws.Range("M89").Value = -6759

$ws.Range("H91").Value = 2979.25
$ws.Range("J91").Value = 2945.875
$ws.Range("L91").Value = 2945.875
$ws.Range("N91").Value = -5753.875

$ws.Range("H92").Value = 1013
$ws.Range("I92").Value = 1149.3334
$ws.Range("J92").Value = 604
$ws.Range("K92").Value = 1149.3334
$ws.Range("L92").Value = 604
$ws.Range("M92").Value = 98.66660000000002
$ws.Range("N92").Value = -3100

$ws.Range("H99").Value = 2455.6
$ws.Range("J99").Value = 3443.8572
$ws.Range("L99").Value = 10331.5716
$ws.Range("N99").Value = -13327.5716

$ws.Range("H102").Value = 57536.5
$ws.Range("J102").Value = 57536.5
$ws.Range("L102").Value = 57536.5
$ws.Range("N102").Value = -64026.5

$ws.Range("H135").Value = 6978.6665
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").Value = ""

$ws.Range("H137").Value = 4350
$ws.Range("I137").Value = 2267.1667
$ws.Range("K137").Value = 6801.500100000001
$ws.Range("M137").Value = -4251.500100000001

$ws.Range("H138").Value = 2357.0476
$ws.Range("I138").Value = 2358.1052
$ws.Range("K138").Value = 7074.3156
$ws.Range("M138").Value = -1934.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3468.3684
$ws.Range("I2").Value = 3676.2307
$ws.Range("J2").Value = 3018
$ws.Range("K2").Value = 3676.2307
$ws.Range("L2").Value = 3018
$ws.Range("M2").Value = -3563.2307
$ws.Range("N2").Value = -3244

$ws.Range("H70").Value = 99999
$ws.Range("J70").Value = 99999
$ws.Range("L70").Value = 99999
$ws.Range("N70").Value = -100539

$ws.Range("H73").Value = 99999
$ws.Range("J73").Value = 99999
$ws.Range("L73").Value = 99999
$ws.Range("N73").Value = -101871

$ws.Range("H116").Value = 3468.3684
$ws.Range("I116").Value = 3676.2307
$ws.Range("J116").Value = 3018
$ws.Range("K116").Value = 3676.2307
$ws.Range("L116").Value = 3018
$ws.Range("M116").Value = -1382.2307
$ws.Range("N116").Value = -7606

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3468.3684
$ws.Range("I3").Value = 3676.2307
$ws.Range("J3").Value = 3018
$ws.Range("K3").Value = 3676.2307
$ws.Range("L3").Value = 3018
$ws.Range("M3").Value = -3562.2307
$ws.Range("N3").Value = -3246

$ws.Range("H107").Value = 5334.3
$ws.Range("I107").Value = 5964.875
$ws.Range("K107").Value = 5964.875
$ws.Range("M107").Value = -4044.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4998.5557
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""

$ws.Range("H31").Value = 100004296
$ws.Range("I31").Value = 125003870
$ws.Range("K31").Value = 125003870
$ws.Range("M31").Value = -125003575

$ws.Range("H34").Value = 100004296
$ws.Range("I34").Value = 125003870
$ws.Range("K34").Value = 125003870
$ws.Range("M34").Value = -125003668

$ws.Range("H132").Value = 2703.5652
$ws.Range("I132").Value = 2484.9048
$ws.Range("K132").Value = 7454.714399999999
$ws.Range("M132").Value = -4924.714399999999

$ws.Range("H134").Value = 2207
$ws.Range("I134").Value = 2207
$ws.Range("K134").Value = 6621
$ws.Range("M134").Value = -4086

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1246.875
$ws.Range("I5").Value = 1062.6666
$ws.Range("K5").Value = 3187.9998
$ws.Range("M5").Value = -3075.9998

$ws.Range("H64").Value = 12364.7
$ws.Range("J64").Value = 13093
$ws.Range("L64").Value = 39279
$ws.Range("N64").Value = -39819

$ws.Range("H67").Value = 12364.7
$ws.Range("J67").Value = 13093
$ws.Range("L67").Value = 39279
$ws.Range("N67").Value = -41151

$ws.Range("H68").Value = 3580.8
$ws.Range("I68").Value = 2451
$ws.Range("J68").Value = 4334
$ws.Range("K68").Value = 7353
$ws.Range("L68").Value = 13002
$ws.Range("M68").Value = -6542
$ws.Range("N68").Value = -14624

$ws.Range("H71").Value = 3580.8
$ws.Range("I71").Value = 2451
$ws.Range("J71").Value = 4334
$ws.Range("K71").Value = 22059
$ws.Range("L71").Value = 39006
$ws.Range("M71").Value = -18003
$ws.Range("N71").Value = -47118

$ws.Range("H112").Value = 15887
$ws.Range("I112").Value = 6438.5
$ws.Range("J112").Value = 19666.4
$ws.Range("K112").Value = 19315.5
$ws.Range("L112").Value = 58999.2
$ws.Range("M112").Value = -18207.5
$ws.Range("N112").Value = -61215.2

$ws.Range("H135").Value = 1246.875
$ws.Range("I135").Value = 1062.6666
$ws.Range("K135").Value = 9563.999400000001
$ws.Range("M135").Value = -7028.999400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""

$ws.Range("H70").Value = 8312.448
$ws.Range("I70").Value = 8210.182000000001
$ws.Range("J70").Value = 8374.944
$ws.Range("K70").Value = 8210.182000000001
$ws.Range("L70").Value = 8374.944
$ws.Range("M70").Value = -7940.182000000001
$ws.Range("N70").Value = -8914.944

$ws.Range("H73").Value = 8312.448
$ws.Range("I73").Value = 8210.182000000001
$ws.Range("K73").Value = 8210.182000000001
$ws.Range("L73").Value = 8374.944
$ws.Range("M73").Value = -7274.182000000001
$ws.Range("N73").Value = -10246.944

$ws.Range("H80").Value = 3455.0908
$ws.Range("I80").Value = 3207.5
$ws.Range("J80").Value = 3596.5715
$ws.Range("K80").Value = 3207.5
$ws.Range("L80").Value = 3596.5715
$ws.Range("M80").Value = -2209.5
$ws.Range("N80").Value = -5592.5715

$ws.Range("H83").Value = 3455.0908
$ws.Range("I83").Value = 3207.5
$ws.Range("J83").Value = 3596.5715
$ws.Range("K83").Value = 16037.5
$ws.Range("L83").Value = 17982.8575
$ws.Range("M83").Value = -11045.5
$ws.Range("N83").Value = -27966.8575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16488.625
$ws.Range("I41").Value = 13982
$ws.Range("J41").Value = 16846.715
$ws.Range("K41").Value = 13982
$ws.Range("L41").Value = 16846.715
$ws.Range("M41").Value = -13592
$ws.Range("N41").Value = -17626.715

$ws.Range("H42").Value = 49974.5
$ws.Range("J42").Value = 49974.5
$ws.Range("L42").Value = 49974.5
$ws.Range("N42").Value = -50730.5

$ws.Range("H132").Value = 479926.2
$ws.Range("J132").Value = 1670231.4
$ws.Range("L132").Value = 5010694.199999999
$ws.Range("N132").Value = -5015754.199999999
